$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7: AssignId, Client Name, Due Date changed
$ws.Range("A7").Value = "assigneng8"
$ws.Range("C7").Value = "wheels india"
$ws.Range("D7").Value = "'2018-05-18"

# Row 8: Client Name, Due Date, Assigner, mode of enquiry changed
$ws.Range("B8").Value = "ds"
$ws.Range("C8").Value = "typical solution"
$ws.Range("D8").Value = "'2018-05-12"
$ws.Range("F8").Value = "Project"
$ws.Range("G8").Value = "wine"

# Reset the quote-prefix styling so the date-like text cells keep the
# default (unstyled) cell format, matching plain text cells elsewhere.
$ws.Range("D7:D8").Style = "Normal"
